# Natmi following Dr Hou advice
# Rebuild the LR-pair data rows: the "sCs" sending-cluster is now included
# as a source cluster (in addition to "FAPs"), giving 6 data rows instead
# of the original 3 (FAPs->ECs, FAPs->M2, FAPs->sCs, sCs->ECs, sCs->M2,
# sCs->sCs), with refreshed statistics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A="FAPs"; B="Tnfsf11"; C="Tnfrsf11a"; D="ECs";
       E=3; F=1; G=1.604474; H=4.813422;
       I=0.9077880937234966; J=0.9077880937234966;
       K=3; L=1; M=1.657081333333333; N=4.971244;
       O=0.07392147563279797; P=0.07392147563279798;
       Q=2.658743915218667; R=23.928695236968;
       S=0.06710503544992556; T=0.06710503544992558 },

    @{ A="FAPs"; B="Tnfsf11"; C="Tnfrsf11a"; D="M2";
       E=3; F=1; G=1.604474; H=4.813422;
       I=0.9077880937234966; J=0.9077880937234966;
       K=3; L=1; M=19.80213966666667; N=59.406419;
       O=0.8833624248860619; P=0.883362424886062;
       Q=31.77201823953533; R=285.948164155818;
       S=0.8019058917542835; T=0.8019058917542836 },

    @{ A="FAPs"; B="Tnfsf11"; C="Tnfrsf11a"; D="sCs";
       E=3; F=1; G=1.604474; H=4.813422;
       I=0.9077880937234966; J=0.9077880937234966;
       K=3; L=1; M=0.9575573333333333; N=2.872672;
       O=0.04271609948114013; P=0.04271609948114013;
       Q=1.536375844842667; R=13.827382603584;
       S=0.03877716651928743; T=0.03877716651928743 },

    @{ A="sCs"; B="Tnfsf11"; C="Tnfrsf11a"; D="ECs";
       E=1; F=0.3333333333333333; G=0.1629803333333333; H=0.488941;
       I=0.09221190627650352; J=0.09221190627650352;
       K=3; L=1; M=1.657081333333333; N=4.971244;
       O=0.07392147563279797; P=0.07392147563279798;
       Q=0.2700716680671111; R=2.430645012604;
       S=0.006816440182872406; T=0.006816440182872406 },

    @{ A="sCs"; B="Tnfsf11"; C="Tnfrsf11a"; D="M2";
       E=1; F=0.3333333333333333; G=0.1629803333333333; H=0.488941;
       I=0.09221190627650352; J=0.09221190627650352;
       K=3; L=1; M=19.80213966666667; N=59.406419;
       O=0.8833624248860619; P=0.883362424886062;
       Q=3.227359323586556; R=29.046233912279;
       S=0.08145653313177842; T=0.08145653313177843 },

    @{ A="sCs"; B="Tnfsf11"; C="Tnfrsf11a"; D="sCs";
       E=1; F=0.3333333333333333; G=0.1629803333333333; H=0.488941;
       I=0.09221190627650352; J=0.09221190627650352;
       K=3; L=1; M=0.9575573333333333; N=2.872672;
       O=0.04271609948114013; P=0.04271609948114013;
       Q=0.1560630133724444; R=1.404567120352;
       S=0.003938932961852695; T=0.003938932961852695 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowData[$c]
    }
}
